$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 101, shifting existing rows 101-111 down to 102-112
$ws.Rows.Item(101).Insert()

# Populate the new row 101 with its data
$ws.Range("A101").Value = 1
$ws.Range("B101").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C101").Value = "Arica y Parinacota"
$ws.Range("D101").Value = 44769
$ws.Range("D101").NumberFormat = $ws.Range("D102").NumberFormat
$ws.Range("E101").Value = 15
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100106
$ws.Range("H101").Value = "Oleaginosos"
$ws.Range("I101").Value = 100106002
$ws.Range("J101").Value = "Palta"
$ws.Range("K101").Value = "Fuerte"
$ws.Range("L101").Value = "Tercera"
$ws.Range("M101").Value = 160
$ws.Range("N101").Value = 44000
$ws.Range("O101").Value = 45000
$ws.Range("P101").Value = 44500
$ws.Range("Q101").Value = "$/caja 25 kilos"
$ws.Range("R101").Value = "Región de Coquimbo"
$ws.Range("S101").Value = 1780
$ws.Range("T101").Value = 25
